# Updated cryptos list (price + 1h volume change refresh, plus a rank
# swap between EthereumClassic and LidoDAOToken at rows 27/28).
#
# Columns D ("Price") and E ("Volume(1h)") are stored as plain text in the
# sheet (values such as "26.937.41" or "1.001" are not valid single
# numbers and must stay text). Assigning a bare numeric-looking string to
# a General-formatted cell makes Excel coerce it to a real number, so
# those assignments are prefixed with a leading apostrophe to force text
# entry, then the cell's Style is reset back to "Normal" so it doesn't
# keep a lingering quote-prefix/text number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.937.41'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.817.21'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'309.63"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = "'0.4685"
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").Value = "'0.07372"
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").Value = "'0.8733"
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").Value = "'20.41"
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '1.811.52'
$ws.Range("E12").Value = '  +4.73%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = "'0.07090"
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = "'6.526"
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("D16").Value = "'91.92"
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '26.974.42'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").Value = "'5.333"
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").Value = "'10.64"
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("D24").Value = '2.038.25'
$ws.Range("E24").Value = '  +3.55%  '
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = "'2.181"
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'18.44"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = "'5.335"
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = "'0.08937"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = "'0.7681"
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = "'1.169"
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("D34").Value = "'4.514"
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("D35").Value = "'2.911"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -2.64%  '
$ws.Range("D38").Value = "'0.01964"
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").Value = "'0.05298"
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("D40").Value = "'2.965"
$ws.Range("E40").Value = '  +2.18%  '
$ws.Range("D41").Value = "'7.279"
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("D42").Value = "'0.5349"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = "'2.330"
$ws.Range("E43").Value = '  -3.46%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = "'8.464"
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").Value = "'0.4935"
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").Value = "'103.40"
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").Value = "'0.06306"
$ws.Range("E51").Value = '  -0.15%  '

# Reset style on cells that required a text-forcing apostrophe prefix
# so they don't retain a distinct quotePrefix style (keeps output clean)
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
